$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert the new columns.
#    Before: A..I  (env2test,host,measurement,app,browserPath,sendMetrics,
#                   telegrafUrl,writeFailReqRspOnly,debug_error)
#    After : A..M  (env2test,host,measurement,app,appComponent,browserPath,
#                   sendMetrics,telegrafUrl,sendResults,xrayUrl,
#                   xrayTestExecKeyMaster,writeFailReqRspOnly,debug_error)
# ---------------------------------------------------------------------------

# New column "appComponent" before the old "browserPath" column (E)
$ws.Columns("E:E").Insert()

# Three new columns "sendResults", "xrayUrl", "xrayTestExecKeyMaster" before
# the old "writeFailReqRspOnly" column (now shifted to I after the first insert)
$ws.Columns("I:K").Insert()

# ---------------------------------------------------------------------------
# 2. Header row (row 1)
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "env2test"
$ws.Range("B1").Value = "host"
$ws.Range("C1").Value = "measurement"
$ws.Range("D1").Value = "app"
$ws.Range("E1").Value = "appComponent"
$ws.Range("F1").Value = "browserPath"
$ws.Range("G1").Value = "sendMetrics"
$ws.Range("H1").Value = "telegrafUrl"
$ws.Range("I1").Value = "sendResults"
$ws.Range("J1").Value = "xrayUrl"
$ws.Range("K1").Value = "xrayTestExecKeyMaster"
$ws.Range("L1").Value = "writeFailReqRspOnly"
$ws.Range("M1").Value = "debug_error"

# ---------------------------------------------------------------------------
# 3. Data row (row 2)
# ---------------------------------------------------------------------------
$ws.Range("A2").Value = "TEST"
$ws.Range("B2").Value = "https://test.clv.cz"
$ws.Range("C2").Value = "test_measurement_mze"
$ws.Range("D2").Value = "test_app"
$ws.Range("E2").Value = "test_comp"
$ws.Range("F2").Value = "C:\Program Files (x86)\Google\Chrome\Application\chrome.exe"
$ws.Range("G2").Value = "false"
$ws.Range("H2").Value = "https://qatick.clance.local/telegraf"
$ws.Range("I2").Value = "false"
$ws.Range("J2").Value = "https://jira.cleverlance.com/rest/raven/1.0/import/execution"
$ws.Range("K2").Value = "n/a"
$ws.Range("L2").Value = "true"
$ws.Range("M2").Value = "false"

# ---------------------------------------------------------------------------
# 4. Cell formatting
#    - header row: gray fill, text format
#    - data row: text format ("@") everywhere, wrap text on the chrome path
# ---------------------------------------------------------------------------
$headerRange = $ws.Range("A1:M1")
$headerRange.NumberFormat = "@"
$headerRange.Interior.ColorIndex = 22

$dataRange = $ws.Range("A2:M2")
$dataRange.NumberFormat = "@"

$ws.Range("F2").WrapText = $true

# ---------------------------------------------------------------------------
# 5. Hyperlinks - rebuild them all so the refs land on the right cells
# ---------------------------------------------------------------------------
$ws.Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("B2"), "https://test.clv.cz/")
$ws.Hyperlinks.Add($ws.Range("H2"), "https://monitoring-test.kb.cz/appmon-in-test/write", [System.Type]::Missing, [System.Type]::Missing, "https://qatick.clance.local/telegraf")
$ws.Hyperlinks.Add($ws.Range("J2"), "https://jira.cleverlance.com/rest/raven/1.0/import/execution")

# ---------------------------------------------------------------------------
# 6. Column widths - let Excel recompute the "best fit" widths like it would
#    have done when the columns were created/edited interactively.
# ---------------------------------------------------------------------------
$ws.Columns("A:M").AutoFit()

# A couple of columns keep an explicit (non bestFit) width in the target file
$ws.Columns("K:K").ColumnWidth = 23.21875

# ---------------------------------------------------------------------------
# 7. Misc sheet bits
# ---------------------------------------------------------------------------
$ws.Range("A2").Select()
